$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 38.04655766666666
$ws.Range("H2").Value = 114.139673
$ws.Range("I2").Value = 0.8090698722086991
$ws.Range("J2").Value = 0.8090698722086992
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 43.69574966666666
$ws.Range("N2").Value = 131.087249
$ws.Range("O2").Value = 0.3365063034544351
$ws.Range("P2").Value = 0.3365063034544351
$ws.Range("Q2").Value = 1662.472859481064
$ws.Range("R2").Value = 14962.25573532958
$ws.Range("S2").Value = 0.2722571119333015
$ws.Range("T2").Value = 0.2722571119333015
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 38.04655766666666
$ws.Range("H3").Value = 114.139673
$ws.Range("I3").Value = 0.8090698722086991
$ws.Range("J3").Value = 0.8090698722086992
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 46.81622333333333
$ws.Range("N3").Value = 140.44867
$ws.Range("O3").Value = 0.3605374521727266
$ws.Range("P3").Value = 0.3605374521727267
$ws.Range("Q3").Value = 1781.196140787212
$ws.Range("R3").Value = 16030.76526708491
$ws.Range("S3").Value = 0.2916999903558379
$ws.Range("T3").Value = 0.291699990355838
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 38.04655766666666
$ws.Range("H4").Value = 114.139673
$ws.Range("I4").Value = 0.8090698722086991
$ws.Range("J4").Value = 0.8090698722086992
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.08903066666667
$ws.Range("N4").Value = 54.26709200000001
$ws.Range("O4").Value = 0.1393058338430899
$ws.Range("P4").Value = 0.1393058338430899
$ws.Range("Q4").Value = 688.2253483934352
$ws.Range("R4").Value = 6194.028135540917
$ws.Range("S4").Value = 0.112708153185355
$ws.Range("T4").Value = 0.112708153185355
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 38.04655766666666
$ws.Range("H5").Value = 114.139673
$ws.Range("I5").Value = 0.8090698722086991
$ws.Range("J5").Value = 0.8090698722086992
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.25020333333334
$ws.Range("N5").Value = 63.75061
$ws.Range("O5").Value = 0.1636504105297484
$ws.Range("P5").Value = 0.1636504105297484
$ws.Range("Q5").Value = 808.497086550059
$ws.Range("R5").Value = 7276.473778950531
$ws.Range("S5").Value = 0.1324046167342046
$ws.Range("T5").Value = 0.1324046167342047
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.617245333333334
$ws.Range("H6").Value = 4.851736000000001
$ws.Range("I6").Value = 0.03439113957782537
$ws.Range("J6").Value = 0.03439113957782537
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 43.69574966666666
$ws.Range("N6").Value = 131.087249
$ws.Range("O6").Value = 0.3365063034544351
$ws.Range("P6").Value = 0.3365063034544351
$ws.Range("Q6").Value = 70.66674723491823
$ws.Range("R6").Value = 636.000725114264
$ws.Range("S6").Value = 0.01157283525091954
$ws.Range("T6").Value = 0.01157283525091954
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.617245333333334
$ws.Range("H7").Value = 4.851736000000001
$ws.Range("I7").Value = 0.03439113957782537
$ws.Range("J7").Value = 0.03439113957782537
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.81622333333333
$ws.Range("N7").Value = 140.44867
$ws.Range("O7").Value = 0.3605374521727266
$ws.Range("P7").Value = 0.3605374521727267
$ws.Range("Q7").Value = 75.71331871012445
$ws.Range("R7").Value = 681.4198683911201
$ws.Range("S7").Value = 0.01239929384070578
$ws.Range("T7").Value = 0.01239929384070578
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.617245333333334
$ws.Range("H8").Value = 4.851736000000001
$ws.Range("I8").Value = 0.03439113957782537
$ws.Range("J8").Value = 0.03439113957782537
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 18.08903066666667
$ws.Range("N8").Value = 54.26709200000001
$ws.Range("O8").Value = 0.1393058338430899
$ws.Range("P8").Value = 0.1393058338430899
$ws.Range("Q8").Value = 29.25440043019023
$ws.Range("R8").Value = 263.2896038717121
$ws.Range("S8").Value = 0.004790886375703053
$ws.Range("T8").Value = 0.004790886375703054
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.617245333333334
$ws.Range("H9").Value = 4.851736000000001
$ws.Range("I9").Value = 0.03439113957782537
$ws.Range("J9").Value = 0.03439113957782537
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 21.25020333333334
$ws.Range("N9").Value = 63.75061
$ws.Range("O9").Value = 0.1636504105297484
$ws.Range("P9").Value = 0.1636504105297484
$ws.Range("Q9").Value = 34.36679217321779
$ws.Range("R9").Value = 309.30112955896
$ws.Range("S9").Value = 0.005628124110496999
$ws.Range("T9").Value = 0.005628124110497
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.361255
$ws.Range("H10").Value = 22.083765
$ws.Range("I10").Value = 0.1565389882134754
$ws.Range("J10").Value = 0.1565389882134754
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 43.69574966666666
$ws.Range("N10").Value = 131.087249
$ws.Range("O10").Value = 0.3365063034544351
$ws.Range("P10").Value = 0.3365063034544351
$ws.Range("Q10").Value = 321.6555557124983
$ws.Range("R10").Value = 2894.900001412484
$ws.Range("S10").Value = 0.052676356270214
$ws.Range("T10").Value = 0.05267635627021401
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 7.361255
$ws.Range("H11").Value = 22.083765
$ws.Range("I11").Value = 0.1565389882134754
$ws.Range("J11").Value = 0.1565389882134754
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 46.81622333333333
$ws.Range("N11").Value = 140.44867
$ws.Range("O11").Value = 0.3605374521727266
$ws.Range("P11").Value = 0.3605374521727267
$ws.Range("Q11").Value = 344.6261580936167
$ws.Range("R11").Value = 3101.63542284255
$ws.Range("S11").Value = 0.05643816797618292
$ws.Range("T11").Value = 0.05643816797618293
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 7.361255
$ws.Range("H12").Value = 22.083765
$ws.Range("I12").Value = 0.1565389882134754
$ws.Range("J12").Value = 0.1565389882134754
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 18.08903066666667
$ws.Range("N12").Value = 54.26709200000001
$ws.Range("O12").Value = 0.1393058338430899
$ws.Range("P12").Value = 0.1393058338430899
$ws.Range("Q12").Value = 133.1579674401534
$ws.Range("R12").Value = 1198.42170696138
$ws.Range("S12").Value = 0.02180679428203182
$ws.Range("T12").Value = 0.02180679428203182
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 7.361255
$ws.Range("H13").Value = 22.083765
$ws.Range("I13").Value = 0.1565389882134754
$ws.Range("J13").Value = 0.1565389882134754
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 21.25020333333334
$ws.Range("N13").Value = 63.75061
$ws.Range("O13").Value = 0.1636504105297484
$ws.Range("P13").Value = 0.1636504105297484
$ws.Range("Q13").Value = 156.4281655385167
$ws.Range("R13").Value = 1407.85348984665
$ws.Range("S13").Value = 0.0256176696850467
$ws.Range("T13").Value = 0.0256176696850467
